# Add 2022-Q4 data:
#  - the existing "2022-Q3" sheet is duplicated (so its data is preserved
#    verbatim as a new "2022-Q3" tab placed right after it)
#  - the original sheet is renamed to "2022-Q4" and repopulated with the new
#    quarter's fund-holding data
#  - the "总计" (totals) sheet gets a new row for 2022-Q4, with the old
#    2022-Q3 totals row shifted down to row 3

$wb = $excel.ActiveWorkbook
$totals = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# 1. Duplicate the current "2022-Q3" sheet right after itself; the copy keeps
#    the old fund data untouched and becomes the new "2022-Q3" tab.
$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item(3)

# 2. Rename the original sheet to "2022-Q4" - it will hold the new data -
#    then rename the copy back to "2022-Q3" now that the name is free.
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

# 3. Clear the old contents of the (now) "2022-Q4" sheet and write the new
#    fund-holding table into it.
$q3.Cells.Clear()

# Copy the header/index cell formatting (bold, centred, bordered) used on the
# "总计" sheet so the new sheet's header row + index column match the target
# styling (style index 2 in the original workbook).
$totals.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$totals.Range("A2").Copy()
$q3.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# B..G are text-typed columns in the target (numeric-looking strings like fund
# codes / percentages must stay literal text, not become real numbers). Mark
# the whole block as Text before writing, then drop back to the default
# "Normal" style afterwards so no stray number-format style sticks to the
# cells (matches the un-styled data cells in the target workbook).
$q3.Range("B2:G9").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "015443"
$q3.Range("C2").Value = "惠升惠享启睿混合A"
$q3.Range("D2").Value = "1.73"
$q3.Range("E2").Value = "64.50"
$q3.Range("F2").Value = "3.65"
$q3.Range("G2").Value = "0.0631"
$q3.Range("H2").Value = 5

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "015110"
$q3.Range("C3").Value = "惠升领先优选混合A"
$q3.Range("D3").Value = "1.49"
$q3.Range("E3").Value = "69.47"
$q3.Range("F3").Value = "3.60"
$q3.Range("G3").Value = "0.0536"
$q3.Range("H3").Value = 7

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "519093"
$q3.Range("C4").Value = "新华钻石品质企业混合"
$q3.Range("D4").Value = "1.34"
$q3.Range("E4").Value = "94.52"
$q3.Range("F4").Value = "3.99"
$q3.Range("G4").Value = "0.0535"
$q3.Range("H4").Value = 9

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "014786"
$q3.Range("C5").Value = "惠升品质优选混合A"
$q3.Range("D5").Value = "1.52"
$q3.Range("E5").Value = "67.63"
$q3.Range("F5").Value = "3.52"
$q3.Range("G5").Value = "0.0535"
$q3.Range("H5").Value = 7

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "001004"
$q3.Range("C6").Value = "新华稳健回报灵活配置混合"
$q3.Range("D6").Value = "0.69"
$q3.Range("E6").Value = "93.81"
$q3.Range("F6").Value = "4.10"
$q3.Range("G6").Value = "0.0283"
$q3.Range("H6").Value = 8

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "015111"
$q3.Range("C7").Value = "惠升领先优选混合C"
$q3.Range("D7").Value = "0.00"
$q3.Range("E7").Value = "69.47"
$q3.Range("F7").Value = "3.60"
$q3.Range("H7").Value = 7

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "015444"
$q3.Range("C8").Value = "惠升惠享启睿混合C"
$q3.Range("D8").Value = "0.00"
$q3.Range("E8").Value = "64.50"
$q3.Range("F8").Value = "3.65"
$q3.Range("H8").Value = 5

$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "014787"
$q3.Range("C9").Value = "惠升品质优选混合C"
$q3.Range("D9").Value = "0.00"
$q3.Range("E9").Value = "67.63"
$q3.Range("F9").Value = "3.52"
$q3.Range("H9").Value = 7

# Drop the temporary Text format back to the default style (keeps the cells'
# stored type as text while leaving no explicit style index behind).
$q3.Range("B2:G9").Style = "Normal"

# G7/G8/G9 are numeric zeros in the target (not text).
$q3.Range("G7").Value = 0
$q3.Range("G8").Value = 0
$q3.Range("G9").Value = 0

# 4. Update the "总计" sheet: insert the new 2022-Q4 total row, pushing the
#    old 2022-Q3 row (A2:D2 = 0, "2022-Q3", 2, 0.25) down to row 3, then
#    overwrite row 2 with the new 2022-Q4 totals.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.25

$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 8

$wb.Worksheets.Item(1).Activate()
